$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.923.90"
$ws.Range("E2").Value = "  -0.61%  "
$ws.Range("D3").Value = "2.545.08"
$ws.Range("E3").Value = "  +3.25%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "567.18"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.02%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.51"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.11%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  -1.17%  "
$ws.Range("D9").Value = "2.544.36"
$ws.Range("E9").Value = "  +3.23%  "
$ws.Range("E10").Value = "  -1.06%  "
$ws.Range("E11").Value = "  -2.25%  "
$ws.Range("E12").Value = "  +0.55%  "
$ws.Range("E13").Value = "  -0.56%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.16"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.67%  "
$ws.Range("D15").Value = "3.000.77"
$ws.Range("E15").Value = "  +3.17%  "
$ws.Range("D16").Value = "62.911.04"
$ws.Range("E16").Value = "  -0.48%  "
$ws.Range("E17").Value = "  -0.83%  "
$ws.Range("D18").Value = "2.539.86"
$ws.Range("E18").Value = "  +2.89%  "
$ws.Range("E19").Value = "  +1.67%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "335.05"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.92%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.29"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.32%  "
$ws.Range("E22").Value = "  -0.77%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.05%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.80"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.23%  "
$ws.Range("E25").Value = "  -3.51%  "
$ws.Range("E26").Value = "  +6.46%  "
$ws.Range("E27").Value = "  +11.54%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.998"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.25%  "
$ws.Range("E29").Value = "  +2.90%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.27"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.96%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.87"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.90%  "
$ws.Range("B32").Value = "PEPE"
$ws.Range("C32").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D32").Value = "0.0₃0807"
$ws.Range("E32").Value = "  -1.10%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "176.49"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.36%  "
$ws.Range("E34").Value = "  +4.48%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "406.54"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +9.24%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.398"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "18.97"
$ws.Range("D37").Style = "Normal"
$ws.Range("E38").Value = "  -0.01%  "
$ws.Range("E39").Value = "  -2.27%  "
$ws.Range("E40").Value = "  +2.61%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.00%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "39.08"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.31%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "152.79"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.92%  "
$ws.Range("E44").Value = "  +0.50%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "20.74"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.57%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.605"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.99%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0955"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.84%  "
$ws.Range("E48").Value = "  -0.62%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0236"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.28%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "18.24"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.08%  "
$ws.Range("E51").Value = "  +0.07%  "
